# "Refresh button + progressbar"
# Adds a new task-log entry (row 8) to the "Feuil1" time sheet: a new
# shared-string task description, a date, and hours spent, then moves the
# active selection down to C9 (mirroring the author's next-row click).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Row 8 was previously blank (sheet only had rows 1-7 then the TOTAL row
# 16). Copy the formatting from row 7 (date / number / wrapped-text styles)
# down into row 8 so the new entry matches the existing table look.
$ws.Range("A7:C7").Copy()
$ws.Range("A8:C8").PasteSpecial()

# New task entry: 2012-06-05, 2 hours, new shared string.
$ws.Range("A8").Value = "6/5/2012"
$ws.Range("B8").Value = 2
$ws.Range("C8").Value = "No pull on refresh but a refresh button with actionView !"

# TOTAL (B16 = SUM(B2:B15)) recalculates automatically to include B8.

# Author's selection moved on to the next row after entering the data.
[void]$ws.Range("C9").Select()
